# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 573-574 (pushing the existing
# rows 573..619 down to 575..621), then populate the two new rows with
# the new observations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows starting at row 573; this shifts rows 573..619
# down to 575..621 and extends the sheet dimension to A1:R621.
$ws.Rows.Item(573).Resize(2).Insert()

# --- New row 573 ---------------------------------------------------
$ws.Cells.Item(573, 1).Value = 5
$ws.Cells.Item(573, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(573, 3).Value = "Maule"
$ws.Cells.Item(573, 4).Value = 45223
$ws.Cells.Item(573, 5).Value = 7
$ws.Cells.Item(573, 6).Value = 100114013
$ws.Cells.Item(573, 7).Value = "Zanahoria"
$ws.Cells.Item(573, 8).Value = "Sin especificar"
$ws.Cells.Item(573, 9).Value = "Primera"
$ws.Cells.Item(573, 10).Value = 400
$ws.Cells.Item(573, 11).Value = 6000
$ws.Cells.Item(573, 12).Value = 6000
$ws.Cells.Item(573, 13).Value = 6000
$ws.Cells.Item(573, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(573, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(573, 16).Value = 300
$ws.Cells.Item(573, 17).Value = 20
$ws.Cells.Item(573, 18).Value = "Hortaliza"

# --- New row 574 ---------------------------------------------------
$ws.Cells.Item(574, 1).Value = 5
$ws.Cells.Item(574, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(574, 3).Value = "Maule"
$ws.Cells.Item(574, 4).Value = 45223
$ws.Cells.Item(574, 5).Value = 7
$ws.Cells.Item(574, 6).Value = 100114013
$ws.Cells.Item(574, 7).Value = "Zanahoria"
$ws.Cells.Item(574, 8).Value = "Sin especificar"
$ws.Cells.Item(574, 9).Value = "Primera"
$ws.Cells.Item(574, 10).Value = 400
$ws.Cells.Item(574, 11).Value = 5000
$ws.Cells.Item(574, 12).Value = 5000
$ws.Cells.Item(574, 13).Value = 5000
$ws.Cells.Item(574, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(574, 15).Value = "Región de Ñuble"
$ws.Cells.Item(574, 16).Value = 250
$ws.Cells.Item(574, 17).Value = 20
$ws.Cells.Item(574, 18).Value = "Hortaliza"
